$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.125.25"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "3.424.29"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.18"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.83"
$ws.Range("E6").Value = "  +5.34%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.07"
$ws.Range("E9").Value = "  +3.78%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.420"
$ws.Range("E11").Value = "  +3.87%  "
$ws.Range("D12").Value = "4.006.86"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.74"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000173"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "3.415.59"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "62.095.34"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.55"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.46"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.00"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.77"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.573"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "76.05"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "3.561.01"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.179"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.69"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.34"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.59"
$ws.Range("E35").Value = "  +5.83%  "
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.98"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.74"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.04"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").Value = "3.458.67"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0785"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.78"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.68"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "2.557.93"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.46"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.82"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.23"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("E51").Value = "  -0.01%  "
